$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.005.26'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '''3.484.69'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("D5").Value = '''577.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").Value = '''161.28'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''3.487.24'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("D9").Value = '''0.581'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.75%  '
$ws.Range("D10").Value = '''7.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").Value = '''0.123'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.91%  '
$ws.Range("E12").Value = '  -1.30%  '
$ws.Range("D13").Value = '''4.084.13'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.81%  '
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '''27.72'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("D16").Value = '''65.091.92'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.0000173'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -9.87%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '''3.481.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("D19").Value = '''6.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.23%  '
$ws.Range("D20").Value = '''13.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.74%  '
$ws.Range("D21").Value = '''382.81'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '''8.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("D23").Value = '''72.92'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '''1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("D25").Value = '''0.535'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.72%  '
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").Value = '''9.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '''1.46'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.39%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''6.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.92%  '
$ws.Range("D32").Value = '''2.03'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.70%  '
$ws.Range("D33").Value = '''23.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").Value = '''7.10'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.35%  '
$ws.Range("E35").Value = '  -1.00%  '
$ws.Range("D36").Value = '''161.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.03%  '
$ws.Range("E37").Value = '  -1.60%  '
$ws.Range("D38").Value = '''0.0757'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.11%  '
$ws.Range("D39").Value = '''27.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("D40").Value = '''2.884.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.10%  '
$ws.Range("D41").Value = '''0.819'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.84%  '
$ws.Range("E42").Value = '  +0.64%  '
$ws.Range("D43").Value = '''4.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("D44").Value = '''43.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("D45").Value = '''26.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("D46").Value = '''0.0311'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.92%  '
$ws.Range("D47").Value = '''2.46'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +12.82%  '
$ws.Range("D48").Value = '''330.33'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.73%  '
$ws.Range("E49").Value = '  -1.47%  '
$ws.Range("D50").Value = '''0.851'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.67%  '
$ws.Range("D51").Value = '''6.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.77%  '
